# Apply updated forecast-error values (rows 2-10) and append a new row 11 (Q9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for existing rows 2-10 (columns B:G)
$data = @{
    2  = @(0.1578884585375505, 1.093026323637141, 3.780296758266214, 1.944298526015543, 1.957160022273444, 51)
    3  = @(0.3634719299439699, 1.11718029823843,  3.396285412598066, 1.842901357261985, 1.825045028446704, 50)
    4  = @(0.2069889977446989, 1.057700982088954, 3.066631775569693, 1.751180109403283, 1.756924272632518, 49)
    5  = @(0.3612381420177023, 1.238224076027852, 3.634105439593412, 1.906332982349467, 1.891601805341011, 48)
    6  = @(0.2411319408316268, 1.211902555076712, 3.480821888571079, 1.865696086872425, 1.870049010965773, 47)
    7  = @(0.3837287414956624, 1.254588007131292, 3.75232404806467,  1.937091646790278, 1.919684637565255, 46)
    8  = @(0.2264048185344497, 1.210514367268617, 3.438743821690994, 1.854385025201345, 1.861309476689316, 45)
    9  = @(0.3609039862410751, 1.204350968061618, 3.365592077467562, 1.834555008024442, 1.81950016883021,  44)
    10 = @(0.2857062726838135, 1.290388063109146, 3.826836026517005, 1.956230054599153, 1.958157187300658, 43)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

# New row 11: label "Q9" with matching style to other A-column labels, plus values
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(11, 1).Value = "Q9"

$ws.Cells.Item(11, 2).Value = 0.3026413976930326
$ws.Cells.Item(11, 3).Value = 1.355495109822078
$ws.Cells.Item(11, 4).Value = 3.884252099312994
$ws.Cells.Item(11, 5).Value = 1.970850602991763
$ws.Cells.Item(11, 6).Value = 1.971081985376132
$ws.Cells.Item(11, 7).Value = 42
